$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, shifting existing rows 140-184 down to 141-185
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new record
$ws.Range("A140").Value = 3
$ws.Range("B140").Value = "Femacal de La Calera"
$ws.Range("C140").Value = "Coquimbo"
$ws.Range("D140").Value = 44524
$ws.Range("E140").Value = 5
$ws.Range("F140").Value = 100112001
$ws.Range("G140").Value = "Berenjena"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 80
$ws.Range("K140").Value = 7500
$ws.Range("L140").Value = 8000
$ws.Range("M140").Value = 7750
$ws.Range("N140").Value = "$/caja 60 unidades"
$ws.Range("O140").Value = "Región de Arica y Parinacota"
$ws.Range("P140").Value = 129
$ws.Range("Q140").Value = 60
$ws.Range("R140").Value = "Hortaliza"
